$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 holds a date serial; bump it by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update price column values on rows 29-32
$ws.Range("D29").Value = 420
$ws.Range("D30").Value = 496
$ws.Range("D31").Value = 872
$ws.Range("D32").Value = 977

$wb.Save()
